# Update Top8_DataComp sheet with refreshed M2/FX length and date figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - China
$ws.Range("C2").Value = 347
$ws.Range("F2").Value = 45597
$ws.Range("G2").Value = 30468
$ws.Range("H2").Value = 45658

# Row 3 - United States
$ws.Range("E3").Value = 30407
$ws.Range("F3").Value = 45597

# Row 4 - Euro Area
$ws.Range("E4").Value = 30407
$ws.Range("F4").Value = 45597
$ws.Range("G4").Value = 30468
$ws.Range("H4").Value = 45658

# Row 5 - Japan
$ws.Range("E5").Value = 30407
$ws.Range("F5").Value = 45597
$ws.Range("G5").Value = 30468
$ws.Range("H5").Value = 45658

# Row 6 - United Kingdom
$ws.Range("C6").Value = 455
$ws.Range("F6").Value = 45566
$ws.Range("G6").Value = 30468
$ws.Range("H6").Value = 45658

# Row 7 - South Korea
$ws.Range("E7").Value = 30376
$ws.Range("F7").Value = 45566
$ws.Range("G7").Value = 30468
$ws.Range("H7").Value = 45658

# Row 8 - Hong Kong
$ws.Range("D8").Value = 423
$ws.Range("E8").Value = 30376
$ws.Range("F8").Value = 45566
$ws.Range("H8").Value = 45658

# Row 9 - Australia
$ws.Range("E9").Value = 30407
$ws.Range("F9").Value = 45597
$ws.Range("G9").Value = 30468
$ws.Range("H9").Value = 45658
